$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.811.77"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "3.373.51"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.372.95"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  -3.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "3.950.89"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "3.372.47"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -4.50%  "
$ws.Range("D18").Value = "60.893.99"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "3.507.83"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.545"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("E27").Value = "  -4.09%  "
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("E29").Value = "  -4.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.67%  "
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  -4.93%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0754"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  -6.48%  "
$ws.Range("D48").Value = "2.537.09"
$ws.Range("E48").Value = "  +8.34%  "
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.26%  "
